# This edit cyclically rotates the content of rows 2-9 on the active sheet
# (the occurrence records), moving each row's data to the next row in the
# cycle: 2 -> 3 -> 9 -> 8 -> 7 -> 6 -> 5 -> 4 -> 2.
#
# Columns that carry per-record data and therefore move with the rotation:
# A, B, D, E, F, G, H, P, Q, R, AN, AO, AW, AX
# (columns C, I, S, T, U, V, W, Y, Z, AA, AB, AD, AE, AG, AT, AY hold the
# same value in every one of these rows, so they are left untouched.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A", "B", "D", "E", "F", "G", "H", "P", "Q", "R", "AN", "AO", "AW", "AX")

# destination row for the content currently sitting in each source row
$destFor = @{ 2 = 3; 3 = 9; 9 = 8; 8 = 7; 7 = 6; 6 = 5; 5 = 4; 4 = 2 }

$sourceRows = @(2, 3, 4, 5, 6, 7, 8, 9)

# 1) Snapshot every source row's current values before any writes happen.
$snapshot = @{}
foreach ($r in $sourceRows) {
    $rowData = @{}
    foreach ($c in $cols) {
        $val = $ws.Range("$c$r").Value2
        $rowData[$c] = $val
    }
    $snapshot[$r] = $rowData
}

# 2) Write each snapshot into its destination row, clearing any cell whose
#    source was empty (so stale data doesn't linger, e.g. AN4/AO4 which
#    move away and are not replaced).
foreach ($r in $sourceRows) {
    $dest = $destFor[$r]
    $rowData = $snapshot[$r]
    foreach ($c in $cols) {
        $val = $rowData[$c]
        $target = $ws.Range("$c$dest")
        if ($null -eq $val -or $val -eq "") {
            $target.ClearContents()
        } else {
            $target.Value = $val
        }
    }
}
